$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop stale outline-level bookkeeping left on the sheet (no rows are
#     actually grouped any more) ---
$ws.Rows("1:8").ClearOutline()

# --- Update tag text on existing "Binary Search" row (row 8, column C) ---
$ws.Range("C8").Value = "#binary-search #必背"

# --- Add new row 9: LeetCode 34 - Find First and Last Position of Element in Sorted Array ---
# Clone formatting from row 8 first so the new row reuses the existing cell
# styles (center alignment / wrap text / date format) instead of creating
# brand-new style records.
$ws.Range("A8:I8").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)

$ws.Range("A9").Value = 34
$ws.Range("B9").Value = "Find First and Last Position of Element in Sorted Array"
$ws.Range("C9").Value = " #array  #binary-search #核心 "
$ws.Range("D9").Value = "medium"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 34
$ws.Range("H9").Value = 44339
$ws.Range("I9").Value = 45831

# --- Row heights: row 8 grows to 34, new row 9 is 51 ---
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 51

# --- Selection moves to F9 ---
$ws.Range("F9").Select()
